# Generate Report for Handback
# Refresh the handback status report timestamps:
#  - Overview sheet: "Latest HO Xliff Generate Date" for the first file
#  - zh-cn / de-de sheets: "Correspond Handoff Datetime" and
#    "Correspond Handback DateTime" for the first file

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-08-21 10:53:44"

$wsZhCn.Range("H2").Value = "2016-08-21 10:53:40"
$wsZhCn.Range("K2").Value = "2016-08-21 10:53:57"

$wsDeDe.Range("H2").Value = "2016-08-21 10:53:44"
$wsDeDe.Range("K2").Value = "2016-08-21 10:54:08"
